$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Inter", 4.83963344788087, 10.4294003868472, 0.6772138458271395, 8, 0.06426506559928842, 76),
    @("Milan", 6.203601108033241, 10.02307692307692, 0.5080711074785451, 16, 0.06399039286326985, 34),
    @("Napoli", 5.401785714285714, 9.677777777777777, 0.5866894562546736, 13, 0.04884758364312267, 10),
    @("Roma", 4.420731707317073, 6.790896159317212, 0.6597321855571497, 11, 0.07802457082418029, 6),
    @("Juventus", 5.12541993281075, 8.488599348534201, 0.6322423256759504, 25, 0.0648176291793313, 43),
    @("Atalanta", 5.121922626025792, 8.354785478547855, 0.5572679897993126, 11, 0.07201680672268908, -1),
    @("Como", 4.064261555806088, 9.308527131782945, 0.6874264032447992, 14, 0.04886548439557035, 15),
    @("Bologna", 4.82559456398641, 8.628037383177571, 0.6164658634538153, 21, 0.1065642070178687, -18),
    @("Lazio", 7.455665024630542, 10.2990990990991, 0.5299703264094956, 23, 0.06404230317273796, -1),
    @("Sassuolo", 5.82172373081464, 7.842756183745583, 0.4127312421004252, 19, 0.08110872455763526, -13),
    @("Udinese", 5.044979079497908, 6.64554794520548, 0.4521527871507626, 18, 0.1002874834144184, -6),
    @("Parma", 5.239157372986369, 7.208487084870849, 0.3830579474051664, 10, 0.08925601507927708, -18),
    @("Cagliari", 4.744791666666667, 9.424063116370808, 0.3594822006472492, 19, 0.083663631494376, -22),
    @("Torino", 5.101123595505618, 7.766208251473477, 0.4010644589000592, 17, 0.09988160585512862, -7),
    @("Cremonese", 5.37816091954023, 8.767206477732794, 0.3702587289625722, 21, 0.09855846235985051, -23),
    @("Genoa", 4.552238805970149, 6.912479740680713, 0.4249416288971295, 12, 0.09199134199134198, -2),
    @("Fiorentina", 4.925790754257908, 9.210922787193974, 0.5047737135771854, 12, 0.07781332329695145, 1),
    @("Lecce", 5.547085201793722, 4.892682926829268, 0.4462581269885185, 23, 0.1044103313840156, -35),
    @("Pisa", 5, 6.701219512195122, 0.3368120263386756, 14, 0.1158613711105489, -17),
    @("Hellas Verona", 4.934545454545455, 5.643527204502814, 0.3868800390911312, 10, 0.09666947216544428, -22)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
